# Actualización automática del mapa (2025-09-11 11:57:00)
#
# Net effect on the "AYKO" sheet:
#   1. A brand-new row is inserted at row 2 (Caso 1497 / Sanchez de
#      Bustamante 2064), pushing the existing data rows down by one.
#   2. The row that used to be case 4768 (VALLESE, FELIPE 684) is removed
#      from the list entirely.
#   3. The row for case 803608463 (Cochabamba /ALT/ 1790) keeps its row
#      position but gets its Estado / Observaciones / Equipo fields updated
#      (node has since been transferred, column can be removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row 2, shifting current rows 2..77 down to 3..78 ---
$ws.Rows.Item(2).Insert()

# --- 2. Remove the old case 4768 row. After the insert above, the data
#        that used to live in row 18 (case 4768) now lives in row 19, so
#        deleting row 19 drops it and shifts everything below back up. ---
$ws.Rows.Item(19).Delete()

# --- 3. Populate the freshly inserted row 2 with its new data ---
# Columns A, B, D and E look numeric/date-like but are stored as plain
# text in this sheet, so a leading apostrophe is used to keep Excel from
# auto-converting them to a number/date.
$ws.Cells.Item(2, 1).Value = "'1497"
$ws.Cells.Item(2, 2).Value = "'4/4/2024"
$ws.Cells.Item(2, 3).Value = "SANCHEZ DE BUSTAMANTE 2064"
$ws.Cells.Item(2, 4).Value = "'2"
$ws.Cells.Item(2, 5).Value = "'784804268"
$ws.Cells.Item(2, 6).Value = "AYKO"
$ws.Cells.Item(2, 7).Value = "Pendiente"
$ws.Cells.Item(2, 8).Value = "Equipo de TLC ya traspasado solo Retirar Columna"
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = "Cambio"
$ws.Cells.Item(2, 11).Value = "Nodo TLC"
$ws.Cells.Item(2, 12).Value = "Pasante"
$ws.Cells.Item(2, 13).Value = -58.406882
$ws.Cells.Item(2, 14).Value = -34.588287
$ws.Cells.Item(2, 15).Value = "Recoleta"
$ws.Cells.Item(2, 16).Value = "Capital Sur"

# --- 4. Update the case 803608463 row (still row 19 after steps 1-2) ---
$ws.Cells.Item(19, 7).Value = "Pendiente"
$ws.Cells.Item(19, 8).Value = "Nodo transferido ya se puede retirar columna"
$ws.Cells.Item(19, 11).Value = "Nodo TLC"
